{"js": "// Find the run text that needs to be split/reworded:\n//   \"t\u1ef1 \u0111\u1ed9ng t\u0103ng n\u1ebfu ki\u1ec3u d\u1eef li\u1ec7u l\u00e0 s\u1ed1\"\n// and turn it into 3 new runs reading:\n//   \"n\u1ebfu \" + \"ki\u1ec3u d\u1eef li\u1ec7u l\u00e0 s\u1ed1\" + \" th\u00ec set t\u1ef1 \u0111\u1ed9ng t\u0103ng\"\n// (the leading \"C\u00e1c field ... kh\u00f3a ch\u00ednh, \" text that precedes it stays a\n// separate run, matching the OOXML diff).\nconst results = context.document.body.search(\n  \"t\u1ef1 \u0111\u1ed9ng t\u0103ng n\u1ebfu ki\u1ec3u d\u1eef li\u1ec7u l\u00e0 s\u1ed1\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found\");\n}\n\nconst target = results.items[0];\n\n// Build a OOXML fragment with three separate <w:r> runs so the replaced\n// text keeps matching run granularity seen in the authored edit.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">n\u1ebfu </w:t></w:r>' +\n  '<w:r><w:t>ki\u1ec3u d\u1eef li\u1ec7u l\u00e0 s\u1ed1</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> th\u00ec set t\u1ef1 \u0111\u1ed9ng t\u0103ng</w:t></w:r>' +\n  '</w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntarget.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "# Locate the run text that needs to be split/reworded:\n#   \"t\u1ef1 \u0111\u1ed9ng t\u0103ng n\u1ebfu ki\u1ec3u d\u1eef li\u1ec7u l\u00e0 s\u1ed1\"\n# and turn it into 3 new runs reading:\n#   \"n\u1ebfu \" + \"ki\u1ec3u d\u1eef li\u1ec7u l\u00e0 s\u1ed1\" + \" th\u00ec set t\u1ef1 \u0111\u1ed9ng t\u0103ng\"\n# (the leading \"C\u00e1c field ... kh\u00f3a ch\u00ednh, \" text that precedes it stays a\n# separate run, matching the OOXML diff).\n$d = $word.ActiveDocument\n\n$finder = $d.Content\n$finder.Find.Text = \"t\u1ef1 \u0111\u1ed9ng t\u0103ng n\u1ebfu ki\u1ec3u d\u1eef li\u1ec7u l\u00e0 s\u1ed1\"\n$found = $finder.Find.Execute()\nif (-not $found) {\n  throw \"Target text not found\"\n}\n\n# Re-materialize a fresh Range object from the found boundaries - Find.Execute\n# repositions $finder in place, but InsertXML needs an independently bound\n# Range to target the right spot.\n$target = $d.Range($finder.Start, $finder.End)\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">n\u1ebfu </w:t></w:r>' +\n  '<w:r><w:t>ki\u1ec3u d\u1eef li\u1ec7u l\u00e0 s\u1ed1</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> th\u00ec set t\u1ef1 \u0111\u1ed9ng t\u0103ng</w:t></w:r>' +\n  '</w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$target.InsertXML($ooxml)\n"}
